# "re do of sceneario 1 2 3 log"
# Updates the "Traceability info" column (D) of the correction grid with the
# re-numbered line references for scenario1.txt / scenario2.txt / scenario3.txt
# (and adds a couple of brand new traceability notes that didn't exist before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- scenario1.txt renumbering (rows in the "Basic Networked Functionality" block) ---
$ws.Range('D21').Value  = 'scenario1.txt(46-60)'
$ws.Range('D23').Value  = 'scenario1.txt(61-65)'
$ws.Range('D32').Value  = 'scenario1.txt(67, 68)'
$ws.Range('D33').Value  = 'scenario1.txt(69)'
$ws.Range('D34').Value  = 'scenario1.txt(70-73)'
$ws.Range('D35').Value  = 'scenario1.txt(74)'
$ws.Range('D36').Value  = 'scenario1.txt(77)'
$ws.Range('D50').Value  = 'scenario1.txt(77)'
$ws.Range('D37').Value  = 'scenario1.txt(82)'
$ws.Range('D38').Value  = 'scenario1.txt(87-89)'
$ws.Range('D40').Value  = 'scenario1.txt(90)'
$ws.Range('D42').Value  = 'scenario1.txt(91)'
$ws.Range('D43').Value  = 'scenario1.txt(92-94)'
$ws.Range('D44').Value  = 'scenario1.txt(95)'
$ws.Range('D45').Value  = 'scenario1.txt(96)'
$ws.Range('D46').Value  = 'scenario1.txt(97-99)'
$ws.Range('D47').Value  = 'scenario1.txt(100)'
$ws.Range('D48').Value  = 'scenario1.txt(101)'
$ws.Range('D51').Value  = 'scenario1.txt(102-111)'
$ws.Range('D52').Value  = 'scenario1.txt(113-123, 118)'
$ws.Range('D53').Value  = 'scenario1.txt(112)'
$ws.Range('D54').Value  = 'scenario1.txt(112)'
$ws.Range('D55').Value  = 'scenario1.txt(127, 128)'
$ws.Range('D56').Value  = 'scenario1.txt(129)'
$ws.Range('D57').Value  = 'scenario1.txt(130)'
$ws.Range('D58').Value  = 'scenario1.txt(132-137)'
$ws.Range('D60').Value  = 'scenario1.txt(143)'
$ws.Range('D61').Value  = 'scenario1.txt(145, 146)'
$ws.Range('D62').Value  = 'scenario1.txt(149)'
$ws.Range('D63').Value  = 'scenario1.txt(150-162)'
$ws.Range('D64').Value  = 'scenario1.txt(163-181)'

# --- scenario2.txt renumbering ---
$ws.Range('D69').Value  = 'scenario2.txt(4-17, 23)'
$ws.Range('D70').Value  = 'scenario2.txt(18-22, 40-43)'
$ws.Range('D71').Value  = 'scenario2.txt(33-37, 49-54)'

# --- scenario3.txt renumbering ---
$ws.Range('D74').Value  = 'scenario3.txt(4, 5)'
$ws.Range('D75').Value  = 'scenario3.txt(7-9)'
$ws.Range('D76').Value  = 'scenario3.txt(10-20)'
$ws.Range('D77').Value  = 'scenario3.txt(21-37)'
$ws.Range('D79').Value  = 'scenario3.txt(38-49)'
$ws.Range('D80').Value  = 'scenario3.txt(59, 60)'
$ws.Range('D81').Value  = 'scenario3.txt(48, 4-56, 58) [only 3 discards out of 5 bids – King Arthur +2 bids]'
$ws.Range('D82').Value  = 'scenario3.txt(50-53)'
$ws.Range('D83').Value  = 'scenario3.txt(62-74, 70)'
$ws.Range('D109').Value = 'scenario3.txt(48-56) [player only discards 3 due to bonus bids]'

# --- JSON bonus row: replace the old placeholder text with real traceability ---
$ws.Range('D131').Value = 'scenario1.txt(76, 85, 103), scenario3.txt(38), merlin_log.txt(90, 95)'

# --- brand new traceability notes (merlin row + Unity/Spring row) ---
$ws.Range('D112').Value = 'scenario1.txt(82-87), merlin_log.txt(80-85), merlin_log(88-92), '
$ws.Range('D130').Value = 'scenario1.txt(15, 31-38, 43), merlin_log(5, 20-33)'

# --- cosmetic: column D got wider to fit the longer traceability strings ---
$ws.Columns.Item(4).ColumnWidth = 93.1

# --- cosmetic: author's cursor/selection ended up on D81 after the edits ---
$ws.Range('D81').Select()
